$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "Total Amount Sold"

# Remove the old data row (row 2) entirely
$ws.Rows("2:2").Delete()
